# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E) and "Correspond Handback
# DateTime" (H) columns on row 2 of each locale sheet (zh-cn, de-de) with
# the freshly generated handback timestamps.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-23 22:50:27"
$ws_zhcn.Range("H2").Value = "2016-03-23 22:50:51"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-23 22:50:31"
$ws_dede.Range("H2").Value = "2016-03-23 22:50:57"
